$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update table description (row 2) ---
$ws.Range("B2").Value = "Tabela responsável por armazenar os dados das entradas e saídas dos veículos no estacionamento"

# --- Update attribute rows (column A) ---
$ws.Range("A5").Value = "id"
$ws.Range("A7").Value = "veiculo_id"
$ws.Range("A8").Value = "entrada"
$ws.Range("A9").Value = "saida"

# --- Update description column (H) for rows 5-7 ---
$ws.Range("H5").Value = "Código de identificador da movimentação"
$ws.Range("H6").Value = "Foreign Key da tabela vagas"
$ws.Range("H7").Value = "Foreign Key da tabela veiculos"

# --- Row 8 (entrada) ---
$ws.Range("C8").Value = "timestamp"
$ws.Range("D8").Value = "sem limite"
$ws.Range("H8").Value = "Data e hora que veículo acessou o estacionamento"

# --- Row 9 (saida) ---
$ws.Range("C9").Value = "timestamp"
$ws.Range("D9").Value = "sem limite"
$ws.Range("H9").Value = "Data e hora que veículo deixou o estacionamento"

# --- Row 10 cleared (valor_pago row removed) ---
$ws.Range("A10:H10").ClearContents() | Out-Null

# --- Rename table header (B1) ---
$ws.Range("B1").Value = "Movimentacoes"

# --- Index table (rows 13-15) ---
$ws.Range("A13").Value = "PRIMARY"
$ws.Range("C13").Value = "Sim"
$ws.Range("D13").Value = "Não"
$ws.Range("E13").Value = "Sim"
$ws.Range("F13").Value = "id"

$ws.Range("A14").Value = "Index_vaga_id"
$ws.Range("C14").Value = "Não"
$ws.Range("D14").Value = "Sim"
$ws.Range("E14").Value = "Não"
$ws.Range("F14").Value = "vaga_id"

$ws.Range("A15").Value = "index_veiculo_id"
$ws.Range("C15").Value = "Não"
$ws.Range("D15").Value = "Sim"
$ws.Range("E15").Value = "Não"
$ws.Range("F15").Value = "veiculo_id"

# Row 15's last edit also dropped the (unused) underline formatting that
# previously decorated F15:H15, so it visually matches rows 13-14.
$ws.Range("F15:H15").Font.Underline = $false

# --- Update selection to match final state ---
$ws.Range("A16:B16").Select() | Out-Null
